$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add new column H header (match the existing header style) ---
$ws.Range("H1").Value = "Full Questions Text"
$ws.Range("H1").Font.Bold = $true
$ws.Range("H1").HorizontalAlignment = -4108
$ws.Range("H1").VerticalAlignment = -4160
$ws.Range("H1").Borders.LineStyle = 1

# --- Data rows ---
# Row 2
$ws.Range("A2").Value = "RTI2025-TRIFED-02"
$ws.Range("B2").Value = "TRIFED (under MoTA)"
$ws.Range("D2").Value = "Tribal enterprises, turnover, ESG investments"
$ws.Range("F2").Value = "Pending"
$ws.Range("G2").Value = "Await response"
$ws.Range("H2").Value = "1. List of Van Dhan Kendras with location nature of products and year of establishment`n2. Annual turnover and profit or loss data of each Kendra`n3. Details of training marketing or value addition support provided to these Kendras`n4. Any record of ESG investment or impact investment in these tribal enterprises`n5. Copies of monitoring or evaluation reports for TRIFED activities in Madhya Pradesh"

# Row 3
$ws.Range("A3").Value = "RTI2025-MoTA-01"
$ws.Range("B3").Value = "Ministry of Tribal Affairs"
$ws.Range("D3").Value = "FRA claims, Van Dhan Kendras, tribal cooperatives"
$ws.Range("F3").Value = "Pending"
$ws.Range("G3").Value = "Await response"
$ws.Range("H3").Value = "1. District-wise number of Individual Forest Rights and Community Forest Rights claims filed approved and rejected`n2. Forest area in hectares recognized under CFR`n3. List of Van Dhan Vikas Kendras including year of establishment number of tribal beneficiaries and revenue or profit or loss details from 2020 to 2024`n4. List of tribal cooperatives supported by the Ministry along with type of products employment generated and annual revenue if available`n5. Copies of evaluation reports or audits related to tribal enterprises supported by the Ministry or TRIFED in Madhya Pradesh"

# Row 4
$ws.Range("A4").Value = "RTI2025-MCA-03"
$ws.Range("B4").Value = "Ministry of Corporate Affairs"
$ws.Range("D4").Value = "CSR expenditures, ESG disclosures in tribal districts"
$ws.Range("F4").Value = "Pending"
$ws.Range("G4").Value = "Await response"
$ws.Range("H4").Value = "1. List of companies reporting CSR expenditure in tribal districts of Madhya Pradesh`n2. Project wise CSR spending with name of implementing agency and district`n3. List of companies that submitted ESG or BRSR disclosures related to tribal areas or forest areas in Madhya Pradesh`n4. Any available ESG impact assessments submitted by these companies"

# Row 5
$ws.Range("A5").Value = "RTI2025-MoEFCC-04"
$ws.Range("B5").Value = "MoEFCC Regional Office Bhopal"
$ws.Range("D5").Value = "MFP revenues, CFR governance, tribal forest livelihoods"
$ws.Range("F5").Value = "Pending"
$ws.Range("G5").Value = "Await response"
$ws.Range("H5").Value = "1. District-wise revenue from sale or trade of minor forest produce including tendu leaves mahua chironji and bamboo`n2. Share of this revenue transferred to tribal gatherers or forest rights holders`n3. List of active Joint Forest Management Committees and CFR Committees under your jurisdiction with available performance data`n4. Copies of evaluation or monitoring reports on tribal livelihoods based on forest resources"

# --- Filed Date column (C): keep as literal text "2025-06-10", not an auto-converted date ---
$ws.Range("C2:C5").NumberFormat = "@"
$ws.Range("C2").Value = "2025-06-10"
$ws.Range("C3").Value = "2025-06-10"
$ws.Range("C4").Value = "2025-06-10"
$ws.Range("C5").Value = "2025-06-10"

# --- Wrap text for the whole used range (header + data) ---
$ws.Range("A1:H5").WrapText = $true

# --- Column widths (best-fit approximations) ---
$ws.Columns.Item(1).ColumnWidth = 17.333333333333332
$ws.Columns.Item(2).ColumnWidth = 26.5
$ws.Columns.Item(3).ColumnWidth = 9.5
$ws.Columns.Item(4).ColumnWidth = 46.0
$ws.Columns.Item(5).ColumnWidth = 12.666666666666666
$ws.Columns.Item(6).ColumnWidth = 6.5
$ws.Columns.Item(7).ColumnWidth = 13.666666666666666
$ws.Columns.Item(8).ColumnWidth = 255.0

# --- Row heights ---
$ws.Rows.Item(2).RowHeight = 72
$ws.Rows.Item(3).RowHeight = 72
$ws.Rows.Item(4).RowHeight = 57.6
$ws.Rows.Item(5).RowHeight = 57.6

# --- Selection (matches end state in the authored file) ---
$ws.Range("B5").Select() | Out-Null
